# Weekly update: insert a new price record as row 541 (shifting existing
# rows 541-575 down to 542-576) on the single worksheet of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 541; everything currently at/after
# row 541 (through 575) shifts down to 542..576.
$ws.Rows.Item(541).Insert()

# Populate the newly inserted row 541 with this week's record.
$ws.Range("A541").Value = 4
$ws.Range("B541").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C541").Value = "Los Lagos"
$ws.Range("D541").Value = 44931
$ws.Range("E541").Value = 10
$ws.Range("F541").Value = 100114001
$ws.Range("G541").Value = "Papa"
$ws.Range("H541").Value = "Patagonia"
$ws.Range("I541").Value = "1a nueva(o)"
$ws.Range("J541").Value = 300
$ws.Range("K541").Value = 13000
$ws.Range("L541").Value = 14000
$ws.Range("M541").Value = 13500
$ws.Range("N541").Value = "$/saco 25 kilos"
$ws.Range("O541").Value = "Región de La Araucanía"
$ws.Range("P541").Value = 540
$ws.Range("Q541").Value = 25
$ws.Range("R541").Value = "Hortaliza"
